# Horarios actualizados Línea 141 - 1181
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173) with the
# latest scrape snapshot taken at 04:03:21.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912" — full refresh: 7 upcoming arrivals (was 4)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:03:21"
$ws1.Range("A3").Value = "Total filas: 7"

$rows1 = @(
    @("04:03:21", "04:46", "215A_EL PATO",  43,  "LP1912"),
    @("04:03:21", "04:53", "11_ETCHEVERRY", 50,  "LP1912"),
    @("04:03:21", "05:16", "17_ROMERO",     73,  "LP1912"),
    @("04:03:21", "05:22", "23_HERNANDEZ",  79,  "LP1912"),
    @("04:03:21", "05:34", "215B_EL PATO",  91,  "LP1912"),
    @("04:03:21", "05:46", "15_ABASTO",     103, "LP1912"),
    @("04:03:21", "05:54", "10_OLMOS",      111, "LP1912")
)

$r = 6
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215" — 2 upcoming arrivals (was 1)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:03:21"
$ws2.Range("A3").Value = "Total filas: 2"

$rows2 = @(
    @("04:03:21", "04:46", "215A_EL PATO", 43, "LP1912"),
    @("04:03:21", "05:34", "215B_EL PATO", 91, "LP1912")
)

$r = 6
foreach ($row in $rows2) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "6203-6173" — gains its first upcoming arrival (was empty, 0 rows)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:03:21"
$ws3.Range("A3").Value = "Total filas: 1"

$ws3.Range("A5").Value = "Hora_Scrap"
$ws3.Range("B5").Value = "Hora_Llegada"
$ws3.Range("C5").Value = "Linea"
$ws3.Range("D5").Value = "Minutos"
$ws3.Range("E5").Value = "Parada"

$ws3.Range("A6").Value = "04:03:21"
$ws3.Range("B6").Value = "05:44"
$ws3.Range("C6").Value = "215A_LA PLATA"
$ws3.Range("D6").Value = 101
$ws3.Range("E6").Value = "L6173"
